$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("LP1912")
$ws.Cells.Item(2, 1).Value = "Última actualización: 13:14:41"
$ws.Cells.Item(3, 1).Value = "Total filas: 199"
$ws.Cells.Item(16, 1).Value = "05:44:02"
$ws.Cells.Item(16, 3).Value = "17X38_ROMERO"
$ws.Cells.Item(16, 4).Value = 56
$ws.Cells.Item(17, 1).Value = "06:38:54"
$ws.Cells.Item(17, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(17, 4).Value = 2
$ws.Cells.Item(41, 1).Value = "06:38:54"
$ws.Cells.Item(41, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(41, 4).Value = 82
$ws.Cells.Item(42, 1).Value = "07:52:32"
$ws.Cells.Item(42, 3).Value = "17_ROMERO"
$ws.Cells.Item(42, 4).Value = 8
$ws.Cells.Item(49, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(50, 3).Value = "15_ABASTO"
$ws.Cells.Item(51, 3).Value = "15_ABASTO"
$ws.Cells.Item(52, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(53, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(54, 3).Value = "10_OLMOS"
$ws.Cells.Item(76, 1).Value = "08:40:59"
$ws.Cells.Item(76, 3).Value = "15X38_ABASTO"
$ws.Cells.Item(76, 4).Value = 37
$ws.Cells.Item(77, 1).Value = "08:30:14"
$ws.Cells.Item(77, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(77, 4).Value = 47
$ws.Cells.Item(78, 1).Value = "08:52:33"
$ws.Cells.Item(78, 3).Value = "15X38_ABASTO"
$ws.Cells.Item(78, 4).Value = 26
$ws.Cells.Item(79, 1).Value = "08:30:14"
$ws.Cells.Item(79, 3).Value = "14_ABASTO"
$ws.Cells.Item(79, 4).Value = 48
$ws.Cells.Item(102, 3).Value = "15_ABASTO"
$ws.Cells.Item(103, 3).Value = "14_ABASTO"
$ws.Cells.Item(118, 3).Value = "15X38_ABASTO"
$ws.Cells.Item(119, 3).Value = "14_ABASTO"
$ws.Cells.Item(147, 1).Value = "11:35:40"
$ws.Cells.Item(147, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(147, 4).Value = 56
$ws.Cells.Item(148, 1).Value = "11:13:01"
$ws.Cells.Item(148, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(148, 4).Value = 78
$ws.Cells.Item(149, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(150, 3).Value = "14_ABASTO"
$ws.Cells.Item(153, 1).Value = "12:33:54"
$ws.Cells.Item(153, 3).Value = "15_ABASTO"
$ws.Cells.Item(153, 4).Value = 1
$ws.Cells.Item(154, 1).Value = "10:56:01"
$ws.Cells.Item(154, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(154, 4).Value = 98
$ws.Cells.Item(155, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(156, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(157, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(158, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(160, 1).Value = "12:47:00"
$ws.Cells.Item(160, 3).Value = "15X38_ABASTO"
$ws.Cells.Item(160, 4).Value = 1
$ws.Cells.Item(161, 1).Value = "12:33:54"
$ws.Cells.Item(161, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(161, 4).Value = 15
$ws.Cells.Item(162, 1).Value = "11:55:01"
$ws.Cells.Item(162, 3).Value = "14_ABASTO"
$ws.Cells.Item(162, 4).Value = 53
$ws.Cells.Item(168, 1).Value = "13:14:41"
$ws.Cells.Item(168, 4).Value = 3
$ws.Cells.Item(169, 1).Value = "13:14:41"
$ws.Cells.Item(169, 4).Value = 5
$ws.Cells.Item(170, 1).Value = "13:14:41"
$ws.Cells.Item(170, 4).Value = 7
$ws.Cells.Item(173, 1).Value = "13:14:41"
$ws.Cells.Item(173, 4).Value = 11
$ws.Cells.Item(174, 1).Value = "13:14:41"
$ws.Cells.Item(174, 3).Value = "14_ABASTO"
$ws.Cells.Item(174, 4).Value = 18
$ws.Cells.Item(175, 1).Value = "12:33:54"
$ws.Cells.Item(175, 3).Value = "215A_EL PATO"
$ws.Cells.Item(175, 4).Value = 59
$ws.Cells.Item(176, 1).Value = "13:14:41"
$ws.Cells.Item(176, 3).Value = "215A_EL PATO"
$ws.Cells.Item(176, 4).Value = 19
$ws.Cells.Item(177, 1).Value = "12:47:00"
$ws.Cells.Item(177, 3).Value = "14_ABASTO"
$ws.Cells.Item(177, 4).Value = 46
$ws.Cells.Item(179, 1).Value = "13:14:41"
$ws.Cells.Item(179, 4).Value = 33
$ws.Cells.Item(181, 1).Value = "13:14:41"
$ws.Cells.Item(181, 4).Value = 40
$ws.Cells.Item(185, 1).Value = "13:14:41"
$ws.Cells.Item(185, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(185, 4).Value = 48
$ws.Cells.Item(186, 1).Value = "13:14:41"
$ws.Cells.Item(186, 2).Value = "14:02"
$ws.Cells.Item(186, 3).Value = "10_OLMOS"
$ws.Cells.Item(186, 4).Value = 48
$ws.Cells.Item(187, 1).Value = "13:14:41"
$ws.Cells.Item(187, 2).Value = "14:02"
$ws.Cells.Item(187, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(187, 4).Value = 48
$ws.Cells.Item(188, 2).Value = "14:06"
$ws.Cells.Item(188, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(188, 4).Value = 72
$ws.Cells.Item(189, 1).Value = "12:47:00"
$ws.Cells.Item(189, 2).Value = "14:07"
$ws.Cells.Item(189, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(189, 4).Value = 80
$ws.Cells.Item(190, 1).Value = "13:14:41"
$ws.Cells.Item(190, 2).Value = "14:16"
$ws.Cells.Item(190, 4).Value = 62
$ws.Cells.Item(191, 2).Value = "14:17"
$ws.Cells.Item(191, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(191, 4).Value = 90
$ws.Cells.Item(192, 1).Value = "13:14:41"
$ws.Cells.Item(192, 2).Value = "14:17"
$ws.Cells.Item(192, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(192, 4).Value = 63
$ws.Cells.Item(193, 1).Value = "12:47:00"
$ws.Cells.Item(193, 2).Value = "14:18"
$ws.Cells.Item(193, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(193, 4).Value = 91
$ws.Cells.Item(194, 1).Value = "13:14:41"
$ws.Cells.Item(194, 2).Value = "14:27"
$ws.Cells.Item(194, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(194, 4).Value = 73
$ws.Cells.Item(195, 1).Value = "12:33:54"
$ws.Cells.Item(195, 2).Value = "14:31"
$ws.Cells.Item(195, 3).Value = "14X44_ABASTO"
$ws.Cells.Item(195, 4).Value = 118
$ws.Cells.Item(196, 1).Value = "13:14:41"
$ws.Cells.Item(196, 2).Value = "14:32"
$ws.Cells.Item(196, 3).Value = "14X44_ABASTO"
$ws.Cells.Item(196, 4).Value = 78
$ws.Cells.Item(197, 1).Value = "13:14:41"
$ws.Cells.Item(197, 2).Value = "14:33"
$ws.Cells.Item(197, 3).Value = "215C_EL PATO"
$ws.Cells.Item(197, 4).Value = 79
$ws.Cells.Item(198, 1).Value = "12:47:00"
$ws.Cells.Item(198, 2).Value = "14:34"
$ws.Cells.Item(198, 3).Value = "215C_EL PATO"
$ws.Cells.Item(198, 4).Value = 107
$ws.Cells.Item(199, 1).Value = "13:14:41"
$ws.Cells.Item(199, 2).Value = "14:39"
$ws.Cells.Item(199, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(199, 4).Value = 85
$ws.Cells.Item(200, 1).Value = "13:14:41"
$ws.Cells.Item(200, 2).Value = "14:47"
$ws.Cells.Item(200, 3).Value = "215B_EL PATO"
$ws.Cells.Item(200, 4).Value = 93
$ws.Cells.Item(200, 5).Value = "LP1912"
$ws.Cells.Item(201, 1).Value = "13:14:41"
$ws.Cells.Item(201, 2).Value = "14:53"
$ws.Cells.Item(201, 3).Value = "215_EL PELIGRO"
$ws.Cells.Item(201, 4).Value = 99
$ws.Cells.Item(201, 5).Value = "LP1912"
$ws.Cells.Item(202, 1).Value = "13:14:41"
$ws.Cells.Item(202, 2).Value = "15:02"
$ws.Cells.Item(202, 3).Value = "10_OLMOS"
$ws.Cells.Item(202, 4).Value = 108
$ws.Cells.Item(202, 5).Value = "LP1912"
$ws.Cells.Item(203, 1).Value = "13:14:41"
$ws.Cells.Item(203, 2).Value = "15:11"
$ws.Cells.Item(203, 3).Value = "14_ABASTO"
$ws.Cells.Item(203, 4).Value = 117
$ws.Cells.Item(203, 5).Value = "LP1912"
$ws.Cells.Item(204, 1).Value = "13:14:41"
$ws.Cells.Item(204, 2).Value = "15:13"
$ws.Cells.Item(204, 3).Value = "17X38_ROMERO"
$ws.Cells.Item(204, 4).Value = 119
$ws.Cells.Item(204, 5).Value = "LP1912"

$ws = $wb.Worksheets.Item("LP1912-215")
$ws.Cells.Item(2, 1).Value = "Última actualización: 13:14:41"
$ws.Cells.Item(30, 1).Value = "13:14:41"
$ws.Cells.Item(30, 4).Value = 19
$ws.Cells.Item(31, 1).Value = "13:14:41"
$ws.Cells.Item(31, 4).Value = 79
$ws.Cells.Item(33, 1).Value = "13:14:41"
$ws.Cells.Item(33, 4).Value = 93
$ws.Cells.Item(34, 1).Value = "13:14:41"
$ws.Cells.Item(34, 4).Value = 99

$ws = $wb.Worksheets.Item("6203-6173")
$ws.Cells.Item(2, 1).Value = "Última actualización: 13:14:41"
$ws.Cells.Item(3, 1).Value = "Total filas: 25"
$ws.Cells.Item(26, 1).Value = "13:14:41"
$ws.Cells.Item(26, 2).Value = "13:16"
$ws.Cells.Item(26, 3).Value = "215C_LA PLATA"
$ws.Cells.Item(26, 4).Value = 2
$ws.Cells.Item(26, 5).Value = "L6203"
$ws.Cells.Item(27, 1).Value = "13:14:41"
$ws.Cells.Item(27, 2).Value = "13:20"
$ws.Cells.Item(27, 4).Value = 6
$ws.Cells.Item(28, 1).Value = "12:47:00"
$ws.Cells.Item(28, 2).Value = "13:21"
$ws.Cells.Item(28, 3).Value = "215B_LP-P MOR-1 Y 57"
$ws.Cells.Item(28, 4).Value = 34
$ws.Cells.Item(28, 5).Value = "L6173"
$ws.Cells.Item(29, 1).Value = "13:14:41"
$ws.Cells.Item(29, 2).Value = "13:56"
$ws.Cells.Item(29, 4).Value = 42
$ws.Cells.Item(30, 1).Value = "12:47:00"
$ws.Cells.Item(30, 2).Value = "13:57"
$ws.Cells.Item(30, 3).Value = "215C_LA PLATA"
$ws.Cells.Item(30, 4).Value = 70
$ws.Cells.Item(30, 5).Value = "L6203"

Write-Host "applied updates"
